$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings stored as text in the source sheet.
# Some of the new values look like plain numbers (e.g. "51.00", "244.04");
# left alone, Excel would auto-convert those to numeric cells and silently
# drop the trailing zero. Force the text format first so they stay text,
# matching the inline-string cells already used throughout the sheet.
$ws.Range("D2").Value = "34.724.22"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "1.863.19"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("E4").Value = "  -0.89%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.04"
$ws.Range("E5").Value = "  -3.96%  "
$ws.Range("E6").Value = "  -7.31%  "
$ws.Range("E7").Value = "  -0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.26"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.00"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0723"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0965"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").Value = "2.137.21"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.67"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.703"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.884.16"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.78"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "34.705.96"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.82"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").Value = "0.0₃0804"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("E23").Value = "  -5.67%  "
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  +4.35%  "
$ws.Range("E26").Value = "  -10.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.79"
$ws.Range("E27").Value = "  -3.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  -3.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.93"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("E30").Value = "  -7.28%  "
$ws.Range("D31").Value = "4.128.44"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.68"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  -5.77%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.811"
$ws.Range("E37").Value = "  -11.66%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("E38").Value = "  -20.83%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "96.32"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0667"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.77"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -6.63%  "
$ws.Range("D45").Value = "1.271.30"
$ws.Range("E45").Value = "  -5.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0808"
$ws.Range("E46").Value = "  +9.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("E47").Value = "  -7.19%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.78"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.22"
